$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update District names to official names
$ws.Range("G5").Value = "Vijayapura (Bijapur)"
$ws.Range("G6").Value = "Chikkamagaluru (Chikmagalur)"
$ws.Range("G8").Value = "Vijayapura (Bijapur)"
$ws.Range("G9").Value = "Vijayapura (Bijapur)"
$ws.Range("G10").Value = "Davangere"
$ws.Range("G11").Value = "Davangere"
$ws.Range("G12").Value = "Vijayapura (Bijapur)"
$ws.Range("G14").Value = "Chikkamagaluru (Chikmagalur)"
$ws.Range("G32").Value = "Davangere"
$ws.Range("G37").Value = "Davangere"
$ws.Range("G40").Value = "Davangere"
$ws.Range("G41").Value = "Vijayapura (Bijapur)"

# Remove stray empty Address (F column) cells that had no content
$ws.Range("F7").ClearContents()
$ws.Range("F15").ClearContents()
$ws.Range("F16").ClearContents()
$ws.Range("F17").ClearContents()
$ws.Range("F31").ClearContents()
$ws.Range("F42").ClearContents()
